# Update cryptos list with latest price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($row, $value) {
    # Force the Price column to remain a text string (matches source data,
    # which is plain text such as "1.00", "60.665.53", etc.) rather than
    # letting Excel auto-convert numeric-looking text into a real number.
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-Volume($row, $value) {
    $ws.Range("E$row").Value = $value
}

Set-PriceText 2  "60.665.53"
Set-Volume    2  "  -1.39%  "

Set-PriceText 3  "2.349.11"
Set-Volume    3  "  -4.75%  "

Set-PriceText 4  "1.00"
Set-Volume    4  "  -0.03%  "

Set-PriceText 5  "542.28"
Set-Volume    5  "  -1.92%  "

Set-PriceText 6  "137.32"
Set-Volume    6  "  -6.64%  "

Set-PriceText 7  "1.00"
Set-Volume    7  "  -0.06%  "

Set-PriceText 8  "0.518"
Set-Volume    8  "  -12.87%  "

Set-PriceText 9  "2.348.83"
Set-Volume    9  "  -4.64%  "

Set-Volume    10 "  -3.99%  "

Set-Volume    11 "  -0.12%  "

Set-PriceText 12 "5.21"
Set-Volume    12 "  -4.93%  "

Set-PriceText 13 "0.338"
Set-Volume    13 "  -4.72%  "

Set-PriceText 14 "24.63"
Set-Volume    14 "  -6.22%  "

Set-PriceText 15 "2.772.85"
Set-Volume    15 "  -4.77%  "

Set-PriceText 16 "60.332.93"
Set-Volume    16 "  -1.81%  "

Set-PriceText 17 "0.0000160"
Set-Volume    17 "  -4.31%  "

Set-PriceText 18 "2.347.55"
Set-Volume    18 "  -4.86%  "

Set-PriceText 19 "10.56"
Set-Volume    19 "  -5.22%  "

Set-PriceText 20 "4.07"
Set-Volume    20 "  -3.38%  "

Set-PriceText 21 "313.29"
Set-Volume    21 "  -2.12%  "

Set-PriceText 22 "6.53"
Set-Volume    22 "  -9.08%  "

Set-PriceText 23 "0.999"
Set-Volume    23 "  -0.22%  "

Set-PriceText 24 "1.86"
Set-Volume    24 "  -1.46%  "

Set-PriceText 25 "62.65"
Set-Volume    25 "  -2.65%  "

Set-PriceText 26 "8.13"
Set-Volume    26 "  +5.33%  "

Set-PriceText 27 "0.999"
Set-Volume    27 "  -0.26%  "

Set-PriceText 28 "2.452.13"
Set-Volume    28 "  -5.92%  "

Set-PriceText 29 "7.92"
Set-Volume    29 "  -4.67%  "

Set-PriceText 30 "0.0₃0885"
Set-Volume    30 "  -11.00%  "

Set-PriceText 31 "502.41"
Set-Volume    31 "  -10.94%  "

Set-Volume    32 "  -8.02%  "

Set-PriceText 33 "0.144"
Set-Volume    33 "  -3.46%  "

Set-PriceText 34 "1.79"
Set-Volume    34 "  -7.16%  "

Set-Volume    35 "  -4.19%  "

Set-PriceText 36 "1.00"
Set-Volume    36 "  +0.06%  "

# Rows 37 and 38 swap coins (NEARProtocol now ranked above EthereumClassic)
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-PriceText 37 "4.52"
Set-Volume    37 "  -7.41%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-PriceText 38 "18.39"
Set-Volume    38 "  -0.38%  "

Set-Volume    39 "  -2.73%  "

Set-Volume    40 "  -11.83%  "

Set-Volume    41 "  +1.00%  "

Set-PriceText 43 "138.24"
Set-Volume    43 "  -3.67%  "

Set-Volume    44 "  -1.17%  "

Set-PriceText 45 "138.21"
Set-Volume    45 "  -6.01%  "

Set-PriceText 46 "3.52"
Set-Volume    46 "  -2.76%  "

Set-PriceText 47 "2.07"
Set-Volume    47 "  -14.38%  "

Set-PriceText 48 "0.0509"
Set-Volume    48 "  -5.49%  "

Set-PriceText 49 "19.40"
Set-Volume    49 "  -10.59%  "

Set-PriceText 50 "0.565"
Set-Volume    50 "  -4.80%  "

Set-PriceText 51 "0.0892"
Set-Volume    51 "  -5.22%  "
